$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (index 1)
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("B2").Value = "2024-06-09"
$ws1.Range("C2").Value = "南昌·第三届龙年动漫展——庆端午贺高考专场"
$ws1.Range("D2").Value = "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆"
$ws1.Range("E2").Value = "2024.06.09 10:00-06.10 18:00"
$ws1.Range("F2").Value = 1678
$ws1.Range("G2").Value = "不可售"
$ws1.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=85297"
$ws1.Range("I2").Value = "//i1.hdslb.com/bfs/openplatform/202405/zBSAcG1V1714936299746.jpeg"

$ws1.Range("B3").Value = "2024-06-15"
$ws1.Range("C3").Value = "上饶·宅舞联萌·随舞动漫派对（免费活动)"
$ws1.Range("D3").Value = "春江北大道和吉阳路交汇处 槠溪时光PARK"
$ws1.Range("E3").Value = "2024.06.15 08:00-06.15 21:00"
$ws1.Range("F3").Value = 34
$ws1.Range("G3").Value = 22.33
$ws1.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=85607"
$ws1.Range("I3").Value = "//i0.hdslb.com/bfs/openplatform/202405/jcZGKqhx1715589649770.jpeg"

$ws1.Range("B4").Value = "2024-06-22"
$ws1.Range("C4").Value = "景德镇·BM次元盛典运动番only"
$ws1.Range("D4").Value = "广场南路金幕影城旁 罗曼园宴会酒店"
$ws1.Range("E4").Value = "2024.06.22 10:00-06.22 17:00"
$ws1.Range("F4").Value = 194
$ws1.Range("G4").Value = 55
$ws1.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=85197"
$ws1.Range("I4").Value = "//i2.hdslb.com/bfs/openplatform/202404/Z6eXz0su1714292081978.png"

$ws1.Range("B5").Value = "2024-06-22"
$ws1.Range("C5").Value = "萍乡·AU9夏至国漫展"
$ws1.Range("D5").Value = "金陵东路18号 萍乡市体育馆"
$ws1.Range("E5").Value = "2024.06.22 10:00-06.22 17:00"
$ws1.Range("F5").Value = 44
$ws1.Range("G5").Value = 45
$ws1.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=86453"
$ws1.Range("I5").Value = "//i1.hdslb.com/bfs/openplatform/202405/iFDRERFO1716547195192.jpeg"

$ws1.Range("B6").Value = "2024-06-23"
$ws1.Range("C6").Value = "上饶·BM次元盛典运动番only"
$ws1.Range("D6").Value = "春江北大道时光PARK内 博悦宴会艺术中心"
$ws1.Range("E6").Value = "2024.06.23 10:00-06.23 17:00"
$ws1.Range("F6").Value = 244
$ws1.Range("G6").Value = 55
$ws1.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=85201"
$ws1.Range("I6").Value = "//i1.hdslb.com/bfs/openplatform/202404/30dgkbjT1714293499693.png"

$ws1.Range("B7").Value = "2024-06-23"
$ws1.Range("C7").Value = "赣州·清风霁月·光夜only"
$ws1.Range("D7").Value = "平安大道 麋鹿LiveHouse"
$ws1.Range("E7").Value = "2024.06.23 14:00-06.23 20:00"
$ws1.Range("F7").Value = 26
$ws1.Range("G7").Value = 158
$ws1.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=86993"
$ws1.Range("I7").Value = "//i1.hdslb.com/bfs/openplatform/202406/PklWR8EP1717429316070.jpeg"

$ws1.Range("B8").Value = "2024-06-29"
$ws1.Range("C8").Value = "南昌·第五人格only"
$ws1.Range("D8").Value = "高处见美好生活公园 百家喜宴高新店"
$ws1.Range("E8").Value = "2024.06.29 10:00-06.29 17:00"
$ws1.Range("F8").Value = 101
$ws1.Range("G8").Value = 65
$ws1.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=87043"
$ws1.Range("I8").Value = "//i0.hdslb.com/bfs/openplatform/202405/zir2PYz81717071721569.jpeg"

$ws1.Range("B9").Value = "2024-06-29"
$ws1.Range("C9").Value = "萍乡·BM次元盛典运动番only"
$ws1.Range("D9").Value = "康庄路3号 萍乡梅园国际大酒店"
$ws1.Range("E9").Value = "2024.06.29 10:00-06.29 17:00"
$ws1.Range("F9").Value = 250
$ws1.Range("G9").Value = 55
$ws1.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=85192"
$ws1.Range("I9").Value = "//i1.hdslb.com/bfs/openplatform/202404/byoupYK21714294780383.png"

$ws1.Range("B10").Value = "2024-06-30"
$ws1.Range("C10").Value = "宜春·BM次元盛典运动番only"
$ws1.Range("D10").Value = "鼓楼西路与官圳路交叉口东120米 地中海宴会酒店(润达店)"
$ws1.Range("E10").Value = "2024.06.30 10:00-06.30 17:00"
$ws1.Range("F10").Value = 242
$ws1.Range("G10").Value = 55
$ws1.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=84636"
$ws1.Range("I10").Value = "//i1.hdslb.com/bfs/openplatform/202405/oaGZXKok1715328213440.png"

$ws1.Range("B11").Value = "2024-07-06"
$ws1.Range("C11").Value = "南昌·次元星球动漫游戏展"
$ws1.Range("D11").Value = "龙蟠街666号融创茂1层 融创茂"
$ws1.Range("E11").Value = "2024.07.06 10:00-07.06 17:00"
$ws1.Range("F11").Value = 12
$ws1.Range("G11").Value = "不可售"
$ws1.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=86405"
$ws1.Range("I11").Value = "//i2.hdslb.com/bfs/openplatform/202405/9ZfGuXJ01716796674559.jpeg"

$ws1.Range("B12").Value = "2024-07-06"
$ws1.Range("C12").Value = "鹰潭·BM次元盛典运动番only"
$ws1.Range("D12").Value = "体育馆东路2号九小隔壁 忆江南•宴会楼"
$ws1.Range("E12").Value = "2024.07.06 10:00-07.06 17:00"
$ws1.Range("F12").Value = 35
$ws1.Range("G12").Value = 55
$ws1.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=85997"
$ws1.Range("I12").Value = "//i1.hdslb.com/bfs/openplatform/202405/4yuR8NQc1716259522268.png"

$ws1.Range("B13").Value = "2024-07-07"
$ws1.Range("C13").Value = "赣州·BM次元盛典运动番only"
$ws1.Range("D13").Value = "米瑞金路2口0号上客天下1楼 上客天下.老虔州"
$ws1.Range("E13").Value = "2024.07.07 10:00-07.07 17:00"
$ws1.Range("F13").Value = 23
$ws1.Range("G13").Value = 55
$ws1.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=86602"
$ws1.Range("I13").Value = "//i1.hdslb.com/bfs/openplatform/202405/Xrq9sfkE1716259438090.png"

$ws1.Range("B14").Value = "2024-07-12"
$ws1.Range("C14").Value = "新余·2024第三届MG动漫嘉年华"
$ws1.Range("D14").Value = "仙女湖大道与五一南路交叉口西约180米 老上海风情街水晶厅"
$ws1.Range("E14").Value = "2024.07.12 10:00-07.13 17:30"
$ws1.Range("F14").Value = 74
$ws1.Range("G14").Value = 55
$ws1.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=86536"
$ws1.Range("I14").Value = "//i0.hdslb.com/bfs/openplatform/202405/11RbfeFq1716813676323.jpeg"

$ws1.Range("B15").Value = "2024-07-13"
$ws1.Range("C15").Value = "南昌·SuperComic动漫游戏博览会"
$ws1.Range("D15").Value = "怀玉山大道1315号 南昌绿地国际博览中心"
$ws1.Range("E15").Value = "2024.07.13 09:00-07.14 17:00"
$ws1.Range("F15").Value = 283
$ws1.Range("G15").Value = 65
$ws1.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=86992"
$ws1.Range("I15").Value = "//i1.hdslb.com/bfs/openplatform/202406/wQTAjelJ1717642148929.jpeg"

$ws1.Range("B16").Value = "2024-07-13"
$ws1.Range("C16").Value = "宜春·COMIC WORLD次元创作同人季·动漫游戏嘉年华"
$ws1.Range("D16").Value = "宜春国际商贸城会展中心 宜春国际商贸城会展中心"
$ws1.Range("E16").Value = "2024.07.13 10:00-07.14 17:00"
$ws1.Range("F16").Value = 37
$ws1.Range("G16").Value = 55
$ws1.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=86667"
$ws1.Range("I16").Value = "//i2.hdslb.com/bfs/openplatform/202405/JEjmQOLw1716737193284.jpeg"

$ws1.Range("B17").Value = "2024-07-14"
$ws1.Range("C17").Value = "吉安·COMIC LIFE次元假日05"
$ws1.Range("D17").Value = "东塘大道与阳明西路交叉路口往西约240米 吉安国际会展中心"
$ws1.Range("E17").Value = "2024.07.14 09:00-07.14 18:00"
$ws1.Range("F17").Value = 450
$ws1.Range("G17").Value = 52.1
$ws1.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=85924"
$ws1.Range("I17").Value = "//i2.hdslb.com/bfs/openplatform/202405/tBNLb2671716182857904.jpeg"

$ws1.Range("B18").Value = "2024-07-19"
$ws1.Range("C18").Value = "赣州·第四届赣州半夏动漫展"
$ws1.Range("D18").Value = "105国道东100米赣州毅德城国际会展中心 赣州毅德城国际会展中心"
$ws1.Range("E18").Value = "2024.07.19 10:00-07.21 17:00"
$ws1.Range("F18").Value = 358
$ws1.Range("G18").Value = 55
$ws1.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=86587"
$ws1.Range("I18").Value = "//i1.hdslb.com/bfs/openplatform/202405/tlfL9oq91717053081587.jpeg"

$ws1.Range("B19").Value = "2024-07-20"
$ws1.Range("C19").Value = "南昌·漫拥动漫嘉年华Pro-追光启航"
$ws1.Range("D19").Value = "小蓝南路420号 洪州体育馆"
$ws1.Range("E19").Value = "2024.07.20 09:00-07.21 17:00"
$ws1.Range("F19").Value = 131
$ws1.Range("G19").Value = 52.5
$ws1.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=85796"
$ws1.Range("I19").Value = "//i1.hdslb.com/bfs/openplatform/202404/FawN3tPD1713364764414.png"

$ws1.Range("B20").Value = "2024-07-21"
$ws1.Range("C20").Value = "乐平·CY境界次元动漫夏时庆"
$ws1.Range("D20").Value = "翥山西路182号 佳佳基大酒店"
$ws1.Range("E20").Value = "2024.07.21 10:00-07.21 17:00"
$ws1.Range("F20").Value = 55
$ws1.Range("G20").Value = 30
$ws1.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=86768"
$ws1.Range("I20").Value = "//i1.hdslb.com/bfs/openplatform/202406/3RWgXosx1717381178470.png"

$ws1.Range("B21").Value = "2024-07-21"
$ws1.Range("C21").Value = "九江·SXD动漫嘉年华"
$ws1.Range("D21").Value = "湓浦街道大中路339号 百嘉洲际酒店"
$ws1.Range("E21").Value = "2024.07.21 10:00-07.21 17:30"
$ws1.Range("F21").Value = 28
$ws1.Range("G21").Value = 45
$ws1.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=86832"
$ws1.Range("I21").Value = "//i2.hdslb.com/bfs/openplatform/202406/Acs2Wqx71717394174913.jpeg"

$ws1.Range("B22").Value = "2024-07-21"
$ws1.Range("C22").Value = "萍乡·NL14动漫游戏展·夏日狂想曲"
$ws1.Range("D22").Value = "公园南路168号(近工行城北分理处) 梅生嘉华酒店"
$ws1.Range("E22").Value = "2024.07.21 10:00-07.21 17:00"
$ws1.Range("F22").Value = 35
$ws1.Range("G22").Value = 40
$ws1.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=86658"
$ws1.Range("I22").Value = "//i1.hdslb.com/bfs/openplatform/202405/bccpK1Zb1716969649865.jpeg"

$ws1.Range("B23").Value = "2024-07-26"
$ws1.Range("C23").Value = "南昌·萌卡动漫展"
$ws1.Range("D23").Value = "八一桥街道青山南路118号蓝海购物广场F1 蓝海展览馆"
$ws1.Range("E23").Value = "2024.07.26 09:00-07.28 17:00"
$ws1.Range("F23").Value = 821
$ws1.Range("G23").Value = 65
$ws1.Range("H23").Value = "https://show.bilibili.com/platform/detail.html?id=86776"
$ws1.Range("I23").Value = "//i0.hdslb.com/bfs/openplatform/202406/WIQIJc741717410349369.jpeg"

$ws1.Range("B24").Value = "2024-07-27"
$ws1.Range("C24").Value = "江西·次元星河动漫游戏嘉年华"
$ws1.Range("D24").Value = "九龙大道1177号 南昌绿地国际博览中心"
$ws1.Range("E24").Value = "2024.07.27 10:00-07.28 17:00"
$ws1.Range("F24").Value = 2614
$ws1.Range("G24").Value = 69
$ws1.Range("H24").Value = "https://show.bilibili.com/platform/detail.html?id=85493"
$ws1.Range("I24").Value = "//i1.hdslb.com/bfs/openplatform/202405/jkKGgOqM1717141906659.png"

$ws1.Range("B25").Value = "2024-07-27"
$ws1.Range("C25").Value = "赣州·马娘only"
$ws1.Range("D25").Value = "火车站广场正对面 赣州友尼宝国际酒店(赣州火车站店)"
$ws1.Range("E25").Value = "2024.07.27 09:00-07.27 17:00"
$ws1.Range("F25").Value = 20
$ws1.Range("G25").Value = 60
$ws1.Range("H25").Value = "https://show.bilibili.com/platform/detail.html?id=86772"
$ws1.Range("I25").Value = "//i0.hdslb.com/bfs/openplatform/202406/BYe9CZzh1717172003064.png"

$ws1.Range("B26").Value = "2024-07-28"
$ws1.Range("C26").Value = "赣州·明日方舟only叙拉古夜宴3.0暨同好交流茶话会"
$ws1.Range("D26").Value = "兴国路恒大帝景西门 江西长庚控股有限公司"
$ws1.Range("E26").Value = "2024.07.28 11:00-07.28 17:00"
$ws1.Range("F26").Value = 54
$ws1.Range("G26").Value = 56
$ws1.Range("H26").Value = "https://show.bilibili.com/platform/detail.html?id=85688"
$ws1.Range("I26").Value = "//i1.hdslb.com/bfs/openplatform/202405/5AFwM8QV1715765287721.png"

$ws1.Range("B27").Value = "2024-08-03"
$ws1.Range("C27").Value = "南昌·幻梦境国际动漫游戏嘉年华1th"
$ws1.Range("D27").Value = "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆"
$ws1.Range("E27").Value = "2024.08.03 09:00-08.04 17:30"
$ws1.Range("F27").Value = 514
$ws1.Range("G27").Value = 64
$ws1.Range("H27").Value = "https://show.bilibili.com/platform/detail.html?id=83980"
$ws1.Range("I27").Value = "//i0.hdslb.com/bfs/openplatform/202403/wRTbRtgD1710755902575.jpeg"

$ws1.Range("B28").Value = "2024-08-03"
$ws1.Range("C28").Value = "景德镇·第十五届瓷都ACG动漫游戏博览会"
$ws1.Range("D28").Value = "迎宾大道与寺山路交叉口东200米 陶博城"
$ws1.Range("E28").Value = "2024.08.03 09:00-08.04 17:00"
$ws1.Range("F28").Value = 846
$ws1.Range("G28").Value = 55
$ws1.Range("H28").Value = "https://show.bilibili.com/platform/detail.html?id=86341"
$ws1.Range("I28").Value = "//i0.hdslb.com/bfs/openplatform/202405/Wd6JiV3I1715953735690.png"

$ws1.Range("B29").Value = "2024-08-03"
$ws1.Range("C29").Value = "景德镇·第十五届瓷都ACG动漫游戏博览会—马正阳内场票"
$ws1.Range("D29").Value = "迎宾大道与寺山路交叉口东200米 陶博城"
$ws1.Range("E29").Value = "2024.08.03 08:30-08.03 17:00"
$ws1.Range("F29").Value = 564
$ws1.Range("G29").Value = "已售罄"
$ws1.Range("H29").Value = "https://show.bilibili.com/platform/detail.html?id=85981"
$ws1.Range("I29").Value = "//i2.hdslb.com/bfs/openplatform/202405/yevI9OGA1716445452947.png"

$ws1.Range("B30").Value = "2024-08-03"
$ws1.Range("C30").Value = "樟树·第二届静卿国风动漫文化展览会"
$ws1.Range("D30").Value = "杏佛路89号 樟树银河国际酒店"
$ws1.Range("E30").Value = "2024.08.03 09:00-08.03 17:00"
$ws1.Range("F30").Value = 446
$ws1.Range("G30").Value = 45
$ws1.Range("H30").Value = "https://show.bilibili.com/platform/detail.html?id=86683"
$ws1.Range("I30").Value = "//i2.hdslb.com/bfs/openplatform/202405/KD1hRj6P1716713054977.jpeg"

$ws1.Range("B31").Value = "2024-08-04"
$ws1.Range("C31").Value = "九江·第一届异次元动漫嘉年华"
$ws1.Range("D31").Value = "长虹西大道兴城广场99号 九江半岛宾馆"
$ws1.Range("E31").Value = "2024.08.04 08:00-08.04 17:00"
$ws1.Range("F31").Value = 251
$ws1.Range("G31").Value = 45
$ws1.Range("H31").Value = "https://show.bilibili.com/platform/detail.html?id=84407"
$ws1.Range("I31").Value = "//i2.hdslb.com/bfs/openplatform/202406/65hJjOfJ1717642614493.jpeg"

$ws1.Range("B32").Value = "2024-08-06"
$ws1.Range("C32").Value = "南昌·第一届异次元动漫嘉年华"
$ws1.Range("D32").Value = "民德路411号 东方豪景花园酒店(民德路店)"
$ws1.Range("E32").Value = "2024.08.06 08:00-08.06 17:00"
$ws1.Range("F32").Value = 378
$ws1.Range("G32").Value = 55
$ws1.Range("H32").Value = "https://show.bilibili.com/platform/detail.html?id=84102"
$ws1.Range("I32").Value = "//i1.hdslb.com/bfs/openplatform/202405/BCA0owUW1716878997961.jpeg"

$ws1.Range("B33").Value = "2024-08-06"
$ws1.Range("C33").Value = "宜春·第三十五届静卿国风动漫文化展览会"
$ws1.Range("D33").Value = "宜阳大道19号(交通银行旁) 宜春安缦文华酒店"
$ws1.Range("E33").Value = "2024.08.06 09:00-08.06 17:00"
$ws1.Range("F33").Value = 446
$ws1.Range("G33").Value = 45
$ws1.Range("H33").Value = "https://show.bilibili.com/platform/detail.html?id=86684"
$ws1.Range("I33").Value = "//i1.hdslb.com/bfs/openplatform/202405/45bGPXfQ1716709212619.jpeg"

$ws1.Range("B34").Value = "2024-08-08"
$ws1.Range("C34").Value = "赣州·第二届异次元动漫嘉年华"
$ws1.Range("D34").Value = "金辉路南3号大坪明德小学体育馆2层东侧201办公室 鲲伍体育·赣州经开区综合体育馆"
$ws1.Range("E34").Value = "2024.08.08 08:00-08.08 17:00"
$ws1.Range("F34").Value = 584
$ws1.Range("G34").Value = 45
$ws1.Range("H34").Value = "https://show.bilibili.com/platform/detail.html?id=84184"
$ws1.Range("I34").Value = "//i1.hdslb.com/bfs/openplatform/202405/ayYIVKwP1716879335847.jpeg"

$ws1.Range("B35").Value = "2024-08-10"
$ws1.Range("C35").Value = "高安·第二届静卿国风动漫文化展览会"
$ws1.Range("D35").Value = "华林中路606号 高安华鼎国际大酒店"
$ws1.Range("E35").Value = "2024.08.10 09:00-08.10 17:00"
$ws1.Range("F35").Value = 421
$ws1.Range("G35").Value = 45
$ws1.Range("H35").Value = "https://show.bilibili.com/platform/detail.html?id=86682"
$ws1.Range("I35").Value = "//i2.hdslb.com/bfs/openplatform/202405/UwvNYGne1716711642772.jpeg"

# Remove the now-superseded trailing rows (36-39) from 展览
$ws1.Range("A36:I39").EntireRow.Delete()

# Sheet 4: 全部类型 (index 4)
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("B2").Value = "2024-06-09"
$ws4.Range("C2").Value = "南昌·第三届龙年动漫展——庆端午贺高考专场"
$ws4.Range("D2").Value = "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆"
$ws4.Range("E2").Value = "2024.06.09 10:00-06.10 18:00"
$ws4.Range("F2").Value = 1678
$ws4.Range("G2").Value = "不可售"
$ws4.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=85297"
$ws4.Range("I2").Value = "//i1.hdslb.com/bfs/openplatform/202405/zBSAcG1V1714936299746.jpeg"

$ws4.Range("B3").Value = "2024-06-15"
$ws4.Range("C3").Value = "上饶·宅舞联萌·随舞动漫派对（免费活动)"
$ws4.Range("D3").Value = "春江北大道和吉阳路交汇处 槠溪时光PARK"
$ws4.Range("E3").Value = "2024.06.15 08:00-06.15 21:00"
$ws4.Range("F3").Value = 34
$ws4.Range("G3").Value = 22.33
$ws4.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=85607"
$ws4.Range("I3").Value = "//i0.hdslb.com/bfs/openplatform/202405/jcZGKqhx1715589649770.jpeg"

$ws4.Range("B4").Value = "2024-06-22"
$ws4.Range("C4").Value = "景德镇·BM次元盛典运动番only"
$ws4.Range("D4").Value = "广场南路金幕影城旁 罗曼园宴会酒店"
$ws4.Range("E4").Value = "2024.06.22 10:00-06.22 17:00"
$ws4.Range("F4").Value = 194
$ws4.Range("G4").Value = 55
$ws4.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=85197"
$ws4.Range("I4").Value = "//i2.hdslb.com/bfs/openplatform/202404/Z6eXz0su1714292081978.png"

$ws4.Range("B5").Value = "2024-06-22"
$ws4.Range("C5").Value = "萍乡·AU9夏至国漫展"
$ws4.Range("D5").Value = "金陵东路18号 萍乡市体育馆"
$ws4.Range("E5").Value = "2024.06.22 10:00-06.22 17:00"
$ws4.Range("F5").Value = 44
$ws4.Range("G5").Value = 45
$ws4.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=86453"
$ws4.Range("I5").Value = "//i1.hdslb.com/bfs/openplatform/202405/iFDRERFO1716547195192.jpeg"

$ws4.Range("B6").Value = "2024-06-23"
$ws4.Range("C6").Value = "上饶·BM次元盛典运动番only"
$ws4.Range("D6").Value = "春江北大道时光PARK内 博悦宴会艺术中心"
$ws4.Range("E6").Value = "2024.06.23 10:00-06.23 17:00"
$ws4.Range("F6").Value = 244
$ws4.Range("G6").Value = 55
$ws4.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=85201"
$ws4.Range("I6").Value = "//i1.hdslb.com/bfs/openplatform/202404/30dgkbjT1714293499693.png"

$ws4.Range("B7").Value = "2024-06-23"
$ws4.Range("C7").Value = "赣州·清风霁月·光夜only"
$ws4.Range("D7").Value = "平安大道 麋鹿LiveHouse"
$ws4.Range("E7").Value = "2024.06.23 14:00-06.23 20:00"
$ws4.Range("F7").Value = 26
$ws4.Range("G7").Value = 158
$ws4.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=86993"
$ws4.Range("I7").Value = "//i1.hdslb.com/bfs/openplatform/202406/PklWR8EP1717429316070.jpeg"

$ws4.Range("B8").Value = "2024-06-29"
$ws4.Range("C8").Value = "南昌·第五人格only"
$ws4.Range("D8").Value = "高处见美好生活公园 百家喜宴高新店"
$ws4.Range("E8").Value = "2024.06.29 10:00-06.29 17:00"
$ws4.Range("F8").Value = 101
$ws4.Range("G8").Value = 65
$ws4.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=87043"
$ws4.Range("I8").Value = "//i0.hdslb.com/bfs/openplatform/202405/zir2PYz81717071721569.jpeg"

$ws4.Range("B9").Value = "2024-06-29"
$ws4.Range("C9").Value = "萍乡·BM次元盛典运动番only"
$ws4.Range("D9").Value = "康庄路3号 萍乡梅园国际大酒店"
$ws4.Range("E9").Value = "2024.06.29 10:00-06.29 17:00"
$ws4.Range("F9").Value = 250
$ws4.Range("G9").Value = 55
$ws4.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=85192"
$ws4.Range("I9").Value = "//i1.hdslb.com/bfs/openplatform/202404/byoupYK21714294780383.png"

$ws4.Range("B10").Value = "2024-06-30"
$ws4.Range("C10").Value = "宜春·BM次元盛典运动番only"
$ws4.Range("D10").Value = "鼓楼西路与官圳路交叉口东120米 地中海宴会酒店(润达店)"
$ws4.Range("E10").Value = "2024.06.30 10:00-06.30 17:00"
$ws4.Range("F10").Value = 242
$ws4.Range("G10").Value = 55
$ws4.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=84636"
$ws4.Range("I10").Value = "//i1.hdslb.com/bfs/openplatform/202405/oaGZXKok1715328213440.png"

$ws4.Range("B11").Value = "2024-07-06"
$ws4.Range("C11").Value = "南昌·次元星球动漫游戏展"
$ws4.Range("D11").Value = "龙蟠街666号融创茂1层 融创茂"
$ws4.Range("E11").Value = "2024.07.06 10:00-07.06 17:00"
$ws4.Range("F11").Value = 12
$ws4.Range("G11").Value = "不可售"
$ws4.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=86405"
$ws4.Range("I11").Value = "//i2.hdslb.com/bfs/openplatform/202405/9ZfGuXJ01716796674559.jpeg"

$ws4.Range("B12").Value = "2024-07-06"
$ws4.Range("C12").Value = "鹰潭·BM次元盛典运动番only"
$ws4.Range("D12").Value = "体育馆东路2号九小隔壁 忆江南•宴会楼"
$ws4.Range("E12").Value = "2024.07.06 10:00-07.06 17:00"
$ws4.Range("F12").Value = 35
$ws4.Range("G12").Value = 55
$ws4.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=85997"
$ws4.Range("I12").Value = "//i1.hdslb.com/bfs/openplatform/202405/4yuR8NQc1716259522268.png"

$ws4.Range("B13").Value = "2024-07-07"
$ws4.Range("C13").Value = "赣州·BM次元盛典运动番only"
$ws4.Range("D13").Value = "米瑞金路2口0号上客天下1楼 上客天下.老虔州"
$ws4.Range("E13").Value = "2024.07.07 10:00-07.07 17:00"
$ws4.Range("F13").Value = 23
$ws4.Range("G13").Value = 55
$ws4.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=86602"
$ws4.Range("I13").Value = "//i1.hdslb.com/bfs/openplatform/202405/Xrq9sfkE1716259438090.png"

$ws4.Range("B14").Value = "2024-07-12"
$ws4.Range("C14").Value = "新余·2024第三届MG动漫嘉年华"
$ws4.Range("D14").Value = "仙女湖大道与五一南路交叉口西约180米 老上海风情街水晶厅"
$ws4.Range("E14").Value = "2024.07.12 10:00-07.13 17:30"
$ws4.Range("F14").Value = 74
$ws4.Range("G14").Value = 55
$ws4.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=86536"
$ws4.Range("I14").Value = "//i0.hdslb.com/bfs/openplatform/202405/11RbfeFq1716813676323.jpeg"

$ws4.Range("B15").Value = "2024-07-13"
$ws4.Range("C15").Value = "南昌·SuperComic动漫游戏博览会"
$ws4.Range("D15").Value = "怀玉山大道1315号 南昌绿地国际博览中心"
$ws4.Range("E15").Value = "2024.07.13 09:00-07.14 17:00"
$ws4.Range("F15").Value = 283
$ws4.Range("G15").Value = 65
$ws4.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=86992"
$ws4.Range("I15").Value = "//i1.hdslb.com/bfs/openplatform/202406/wQTAjelJ1717642148929.jpeg"

$ws4.Range("B16").Value = "2024-07-13"
$ws4.Range("C16").Value = "宜春·COMIC WORLD次元创作同人季·动漫游戏嘉年华"
$ws4.Range("D16").Value = "宜春国际商贸城会展中心 宜春国际商贸城会展中心"
$ws4.Range("E16").Value = "2024.07.13 10:00-07.14 17:00"
$ws4.Range("F16").Value = 37
$ws4.Range("G16").Value = 55
$ws4.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=86667"
$ws4.Range("I16").Value = "//i2.hdslb.com/bfs/openplatform/202405/JEjmQOLw1716737193284.jpeg"

$ws4.Range("B17").Value = "2024-07-14"
$ws4.Range("C17").Value = "吉安·COMIC LIFE次元假日05"
$ws4.Range("D17").Value = "东塘大道与阳明西路交叉路口往西约240米 吉安国际会展中心"
$ws4.Range("E17").Value = "2024.07.14 09:00-07.14 18:00"
$ws4.Range("F17").Value = 450
$ws4.Range("G17").Value = 52.1
$ws4.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=85924"
$ws4.Range("I17").Value = "//i2.hdslb.com/bfs/openplatform/202405/tBNLb2671716182857904.jpeg"

$ws4.Range("B18").Value = "2024-07-19"
$ws4.Range("C18").Value = "赣州·第四届赣州半夏动漫展"
$ws4.Range("D18").Value = "105国道东100米赣州毅德城国际会展中心 赣州毅德城国际会展中心"
$ws4.Range("E18").Value = "2024.07.19 10:00-07.21 17:00"
$ws4.Range("F18").Value = 359
$ws4.Range("G18").Value = 55
$ws4.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=86587"
$ws4.Range("I18").Value = "//i1.hdslb.com/bfs/openplatform/202405/tlfL9oq91717053081587.jpeg"

$ws4.Range("B19").Value = "2024-07-20"
$ws4.Range("C19").Value = "南昌·漫拥动漫嘉年华Pro-追光启航"
$ws4.Range("D19").Value = "小蓝南路420号 洪州体育馆"
$ws4.Range("E19").Value = "2024.07.20 09:00-07.21 17:00"
$ws4.Range("F19").Value = 132
$ws4.Range("G19").Value = 52.5
$ws4.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=85796"
$ws4.Range("I19").Value = "//i1.hdslb.com/bfs/openplatform/202404/FawN3tPD1713364764414.png"

$ws4.Range("B20").Value = "2024-07-21"
$ws4.Range("C20").Value = "乐平·CY境界次元动漫夏时庆"
$ws4.Range("D20").Value = "翥山西路182号 佳佳基大酒店"
$ws4.Range("E20").Value = "2024.07.21 10:00-07.21 17:00"
$ws4.Range("F20").Value = 55
$ws4.Range("G20").Value = 30
$ws4.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=86768"
$ws4.Range("I20").Value = "//i1.hdslb.com/bfs/openplatform/202406/3RWgXosx1717381178470.png"

$ws4.Range("B21").Value = "2024-07-21"
$ws4.Range("C21").Value = "九江·SXD动漫嘉年华"
$ws4.Range("D21").Value = "湓浦街道大中路339号 百嘉洲际酒店"
$ws4.Range("E21").Value = "2024.07.21 10:00-07.21 17:30"
$ws4.Range("F21").Value = 28
$ws4.Range("G21").Value = 45
$ws4.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=86832"
$ws4.Range("I21").Value = "//i2.hdslb.com/bfs/openplatform/202406/Acs2Wqx71717394174913.jpeg"

$ws4.Range("B22").Value = "2024-07-21"
$ws4.Range("C22").Value = "萍乡·NL14动漫游戏展·夏日狂想曲"
$ws4.Range("D22").Value = "公园南路168号(近工行城北分理处) 梅生嘉华酒店"
$ws4.Range("E22").Value = "2024.07.21 10:00-07.21 17:00"
$ws4.Range("F22").Value = 35
$ws4.Range("G22").Value = 40
$ws4.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=86658"
$ws4.Range("I22").Value = "//i1.hdslb.com/bfs/openplatform/202405/bccpK1Zb1716969649865.jpeg"

$ws4.Range("B23").Value = "2024-07-26"
$ws4.Range("C23").Value = "南昌·萌卡动漫展"
$ws4.Range("D23").Value = "八一桥街道青山南路118号蓝海购物广场F1 蓝海展览馆"
$ws4.Range("E23").Value = "2024.07.26 09:00-07.28 17:00"
$ws4.Range("F23").Value = 821
$ws4.Range("G23").Value = 65
$ws4.Range("H23").Value = "https://show.bilibili.com/platform/detail.html?id=86776"
$ws4.Range("I23").Value = "//i0.hdslb.com/bfs/openplatform/202406/WIQIJc741717410349369.jpeg"

$ws4.Range("B24").Value = "2024-07-27"
$ws4.Range("C24").Value = "江西·次元星河动漫游戏嘉年华"
$ws4.Range("D24").Value = "九龙大道1177号 南昌绿地国际博览中心"
$ws4.Range("E24").Value = "2024.07.27 10:00-07.28 17:00"
$ws4.Range("F24").Value = 2614
$ws4.Range("G24").Value = 69
$ws4.Range("H24").Value = "https://show.bilibili.com/platform/detail.html?id=85493"
$ws4.Range("I24").Value = "//i1.hdslb.com/bfs/openplatform/202405/jkKGgOqM1717141906659.png"

$ws4.Range("B25").Value = "2024-07-27"
$ws4.Range("C25").Value = "赣州·马娘only"
$ws4.Range("D25").Value = "火车站广场正对面 赣州友尼宝国际酒店(赣州火车站店)"
$ws4.Range("E25").Value = "2024.07.27 09:00-07.27 17:00"
$ws4.Range("F25").Value = 20
$ws4.Range("G25").Value = 60
$ws4.Range("H25").Value = "https://show.bilibili.com/platform/detail.html?id=86772"
$ws4.Range("I25").Value = "//i0.hdslb.com/bfs/openplatform/202406/BYe9CZzh1717172003064.png"

$ws4.Range("B26").Value = "2024-07-28"
$ws4.Range("C26").Value = "赣州·明日方舟only叙拉古夜宴3.0暨同好交流茶话会"
$ws4.Range("D26").Value = "兴国路恒大帝景西门 江西长庚控股有限公司"
$ws4.Range("E26").Value = "2024.07.28 11:00-07.28 17:00"
$ws4.Range("F26").Value = 54
$ws4.Range("G26").Value = 56
$ws4.Range("H26").Value = "https://show.bilibili.com/platform/detail.html?id=85688"
$ws4.Range("I26").Value = "//i1.hdslb.com/bfs/openplatform/202405/5AFwM8QV1715765287721.png"

$ws4.Range("B27").Value = "2024-08-03"
$ws4.Range("C27").Value = "南昌·幻梦境国际动漫游戏嘉年华1th"
$ws4.Range("D27").Value = "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆"
$ws4.Range("E27").Value = "2024.08.03 09:00-08.04 17:30"
$ws4.Range("F27").Value = 514
$ws4.Range("G27").Value = 64
$ws4.Range("H27").Value = "https://show.bilibili.com/platform/detail.html?id=83980"
$ws4.Range("I27").Value = "//i0.hdslb.com/bfs/openplatform/202403/wRTbRtgD1710755902575.jpeg"

$ws4.Range("B28").Value = "2024-08-03"
$ws4.Range("C28").Value = "景德镇·第十五届瓷都ACG动漫游戏博览会"
$ws4.Range("D28").Value = "迎宾大道与寺山路交叉口东200米 陶博城"
$ws4.Range("E28").Value = "2024.08.03 09:00-08.04 17:00"
$ws4.Range("F28").Value = 846
$ws4.Range("G28").Value = 55
$ws4.Range("H28").Value = "https://show.bilibili.com/platform/detail.html?id=86341"
$ws4.Range("I28").Value = "//i0.hdslb.com/bfs/openplatform/202405/Wd6JiV3I1715953735690.png"

$ws4.Range("B29").Value = "2024-08-03"
$ws4.Range("C29").Value = "景德镇·第十五届瓷都ACG动漫游戏博览会—马正阳内场票"
$ws4.Range("D29").Value = "迎宾大道与寺山路交叉口东200米 陶博城"
$ws4.Range("E29").Value = "2024.08.03 08:30-08.03 17:00"
$ws4.Range("F29").Value = 564
$ws4.Range("G29").Value = "已售罄"
$ws4.Range("H29").Value = "https://show.bilibili.com/platform/detail.html?id=85981"
$ws4.Range("I29").Value = "//i2.hdslb.com/bfs/openplatform/202405/yevI9OGA1716445452947.png"

$ws4.Range("B30").Value = "2024-08-03"
$ws4.Range("C30").Value = "樟树·第二届静卿国风动漫文化展览会"
$ws4.Range("D30").Value = "杏佛路89号 樟树银河国际酒店"
$ws4.Range("E30").Value = "2024.08.03 09:00-08.03 17:00"
$ws4.Range("F30").Value = 446
$ws4.Range("G30").Value = 45
$ws4.Range("H30").Value = "https://show.bilibili.com/platform/detail.html?id=86683"
$ws4.Range("I30").Value = "//i2.hdslb.com/bfs/openplatform/202405/KD1hRj6P1716713054977.jpeg"

$ws4.Range("B31").Value = "2024-08-04"
$ws4.Range("C31").Value = "九江·第一届异次元动漫嘉年华"
$ws4.Range("D31").Value = "长虹西大道兴城广场99号 九江半岛宾馆"
$ws4.Range("E31").Value = "2024.08.04 08:00-08.04 17:00"
$ws4.Range("F31").Value = 251
$ws4.Range("G31").Value = 45
$ws4.Range("H31").Value = "https://show.bilibili.com/platform/detail.html?id=84407"
$ws4.Range("I31").Value = "//i2.hdslb.com/bfs/openplatform/202406/65hJjOfJ1717642614493.jpeg"

$ws4.Range("B32").Value = "2024-08-06"
$ws4.Range("C32").Value = "南昌·第一届异次元动漫嘉年华"
$ws4.Range("D32").Value = "民德路411号 东方豪景花园酒店(民德路店)"
$ws4.Range("E32").Value = "2024.08.06 08:00-08.06 17:00"
$ws4.Range("F32").Value = 378
$ws4.Range("G32").Value = 55
$ws4.Range("H32").Value = "https://show.bilibili.com/platform/detail.html?id=84102"
$ws4.Range("I32").Value = "//i1.hdslb.com/bfs/openplatform/202405/BCA0owUW1716878997961.jpeg"

$ws4.Range("B33").Value = "2024-08-06"
$ws4.Range("C33").Value = "宜春·第三十五届静卿国风动漫文化展览会"
$ws4.Range("D33").Value = "宜阳大道19号(交通银行旁) 宜春安缦文华酒店"
$ws4.Range("E33").Value = "2024.08.06 09:00-08.06 17:00"
$ws4.Range("F33").Value = 446
$ws4.Range("G33").Value = 45
$ws4.Range("H33").Value = "https://show.bilibili.com/platform/detail.html?id=86684"
$ws4.Range("I33").Value = "//i1.hdslb.com/bfs/openplatform/202405/45bGPXfQ1716709212619.jpeg"

$ws4.Range("B34").Value = "2024-08-08"
$ws4.Range("C34").Value = "赣州·第二届异次元动漫嘉年华"
$ws4.Range("D34").Value = "金辉路南3号大坪明德小学体育馆2层东侧201办公室 鲲伍体育·赣州经开区综合体育馆"
$ws4.Range("E34").Value = "2024.08.08 08:00-08.08 17:00"
$ws4.Range("F34").Value = 584
$ws4.Range("G34").Value = 45
$ws4.Range("H34").Value = "https://show.bilibili.com/platform/detail.html?id=84184"
$ws4.Range("I34").Value = "//i1.hdslb.com/bfs/openplatform/202405/ayYIVKwP1716879335847.jpeg"

$ws4.Range("B35").Value = "2024-08-10"
$ws4.Range("C35").Value = "高安·第二届静卿国风动漫文化展览会"
$ws4.Range("D35").Value = "华林中路606号 高安华鼎国际大酒店"
$ws4.Range("E35").Value = "2024.08.10 09:00-08.10 17:00"
$ws4.Range("F35").Value = 421
$ws4.Range("G35").Value = 45
$ws4.Range("H35").Value = "https://show.bilibili.com/platform/detail.html?id=86682"
$ws4.Range("I35").Value = "//i2.hdslb.com/bfs/openplatform/202405/UwvNYGne1716711642772.jpeg"

# Remove the now-superseded trailing rows (36-39) from 全部类型
$ws4.Range("A36:I39").EntireRow.Delete()

Write-Output "done"
